$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductData")

# --- Update row 11: product renamed, quantity/prices become real numbers,
#     creation_date bumped to the re-save timestamp ---
$ws.Range("C11").Value = "Resma x 200 hojas"
$ws.Range("D11").Value = 10
$ws.Range("E11").Value = 10000
$ws.Range("F11").Value = 15000
$ws.Range("G11").Value = 45818.97881371528

# --- New row 12: Yagurt 1.5 litro ---
$ws.Range("A12").Value = "6AYB"
$ws.Range("B12").Value = "Alimentos y bebidas"
$ws.Range("C12").Value = "Yagurt 1.5 litro"
$ws.Range("D12").Value = 15
$ws.Range("E12").Value = 1400
$ws.Range("F12").Value = 2500
$ws.Range("G12").Value = 45818.98035993055
$ws.Range("G12").NumberFormat = $ws.Range("G11").NumberFormat

# --- New row 13: Micropunta 0.6mm Good lines (quantity/prices kept as text,
#     matching the source data for this row) ---
$ws.Range("A13").Value = "3P"
$ws.Range("B13").Value = "Papelería"
$ws.Range("C13").Value = "Micropunta 0.6mm Good lines"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1200"
$ws.Range("E13").Style = "Normal"

$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "2400"
$ws.Range("F13").Style = "Normal"

$ws.Range("G13").Value = 45818.98101949068
$ws.Range("G13").NumberFormat = $ws.Range("G11").NumberFormat
